$d = $word.ActiveDocument

# ---------------------------------------------------------------------
# 1) Version-history table: bump the last-revision date
#    09/07/2020 -> 21/09/2020
# ---------------------------------------------------------------------
$d.Content.Find.Execute("09/07/2020", $true, $false, $false, $false, $false, `
                         $true, 1, $false, "21/09/2020", 2) | Out-Null

# ---------------------------------------------------------------------
# 2) Trim the alternate-flow cross-reference on step 13
#    af[3,5] -> af[3]
# ---------------------------------------------------------------------
$d.Content.Find.Execute("af[3,5]", $true, $false, $false, $false, $false, `
                         $true, 1, $false, "af[3]", 2) | Out-Null

# ---------------------------------------------------------------------
# 3) Remove the whole AF[5] "Periodo com dias alternados" alternate flow:
#      - the blank paragraph right before its heading
#      - the "AF[5] - ..." heading paragraph
#      - the "1. Chefe/Beneficiario ... ALTERNADAS ..." body paragraph
#    (the two trailing blank paragraphs before "Exception Flows" stay)
# ---------------------------------------------------------------------
$af5 = $null
$n = $d.Paragraphs.Count
for ($i = 1; $i -le $n; $i++) {
    $p = $d.Paragraphs.Item($i)
    if ($p.Range.Text.Contains("AF[5]")) {
        $af5 = $i
        break
    }
}

if ($af5 -ne $null) {
    $startPara = $d.Paragraphs.Item($af5 - 1)
    $endPara = $d.Paragraphs.Item($af5 + 1)
    $r = $d.Range($startPara.Range.Start, $endPara.Range.End)
    $r.Delete()
}
